$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.313.27'
$ws.Range("E2").Value = '  +3.81%  '

$ws.Range("D3").Value = '2.318.16'
$ws.Range("E3").Value = '  +2.30%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '516.72'
$ws.Range("E5").Value = '  +3.92%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.06'
$ws.Range("E6").Value = '  +3.31%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.18%  '

$ws.Range("E8").Value = '  +2.09%  '

$ws.Range("D9").Value = '2.341.19'
$ws.Range("E9").Value = '  +3.19%  '

$ws.Range("E10").Value = '  +8.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.154'
$ws.Range("E11").Value = '  +1.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.15'
$ws.Range("E12").Value = '  +7.99%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.343'
$ws.Range("E13").Value = '  +2.16%  '

$ws.Range("E14").Value = '  +5.36%  '

$ws.Range("D15").Value = '2.752.83'
$ws.Range("E15").Value = '  +3.25%  '

$ws.Range("D16").Value = '56.378.26'
$ws.Range("E16").Value = '  +3.99%  '

$ws.Range("D18").Value = '2.337.81'
$ws.Range("E18").Value = '  +3.06%  '

$ws.Range("E19").Value = '  +2.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.31'
$ws.Range("E21").Value = '  +6.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.69'
$ws.Range("E22").Value = '  +6.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.93'
$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.993'
$ws.Range("E25").Value = '  -0.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.157'
$ws.Range("E26").Value = '  +5.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.64'
$ws.Range("E27").Value = '  +4.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '172.04'
$ws.Range("E28").Value = '  +0.63%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.19'
$ws.Range("E29").Value = '  +9.77%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.68'
$ws.Range("E30").Value = '  +4.51%  '

$ws.Range("D31").Value = '0.0₃0722'
$ws.Range("E31").Value = '  +4.37%  '

$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.24'
$ws.Range("E32").Value = '  +5.24%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.35'
$ws.Range("E33").Value = '  +3.39%  '

$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.994'
$ws.Range("E35").Value = '  -0.23%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.26'
$ws.Range("E36").Value = '  +5.65%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.930'
$ws.Range("E37").Value = '  -1.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.99'
$ws.Range("E38").Value = '  +7.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.51'
$ws.Range("E39").Value = '  +8.47%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.37'
$ws.Range("E40").Value = '  +3.89%  '

$ws.Range("E41").Value = '  +2.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '138.70'
$ws.Range("E43").Value = '  +10.47%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.06'
$ws.Range("E44").Value = '  +5.45%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '267.77'
$ws.Range("E45").Value = '  +10.99%  '

$ws.Range("E46").Value = '  +2.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0924'
$ws.Range("E47").Value = '  +3.76%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.555'
$ws.Range("E48").Value = '  +1.53%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.383'
$ws.Range("E49").Value = '  +2.69%  '

$ws.Range("E50").Value = '  +5.30%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.86'
$ws.Range("E51").Value = '  +4.48%  '
